# Reorders the 17 data rows (rows 2-18) of the sheet according to the
# permutation observed between the "before" and "after" workbook states.
# Row 1 (header) and columns A,B,C,E,F,G,O,R (constant across all data
# rows) are left untouched; only the per-row data is shuffled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 18
$lastCol = 18   # column R

# Destination row (index) -> source row (value taken from)
$order = @{
    2  = 3
    3  = 11
    4  = 12
    5  = 7
    6  = 13
    7  = 18
    8  = 16
    9  = 17
    10 = 6
    11 = 2
    12 = 4
    13 = 8
    14 = 10
    15 = 14
    16 = 5
    17 = 9
    18 = 15
}

# Snapshot all existing data rows (values only) before overwriting anything,
# since the reorder is a permutation and destinations/sources overlap.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

# Write back each destination row using the value captured from its source row.
foreach ($destRow in $order.Keys) {
    $srcRow = $order[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c - 1]
    }
}
